$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the paragraph index (1-based, Word COM style) whose visible
# text matches a distinctive substring. Re-resolved every time it's needed
# since earlier edits can shift paragraph indices.
# ---------------------------------------------------------------------------
function Get-ParaIndex($needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        if ($d.Paragraphs($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) "If we pass more width than parents width, ... width() function"
#    -> bold the whole paragraph (pPr/rPr + every run's rPr)
# ---------------------------------------------------------------------------
$idx = Get-ParaIndex "If we pass more width than parents width"
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">If we pass more width than parents width, automatically set to parent width for </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>width(</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>) function</w:t>
  </w:r>
</w:p>
'@
$null = $d.Paragraphs($idx).Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 2) ".requiredWidth will actually taken according to given width, ..."
#    -> bold the whole paragraph (pPr/rPr + every run's rPr)
# ---------------------------------------------------------------------------
$idx = Get-ParaIndex "requiredWidth will actually taken"
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>requiredWidth</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve"> will actually taken according to given width, so it does not rely on parents width.</w:t>
  </w:r>
</w:p>
'@
$null = $d.Paragraphs($idx).Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 3) "2 paremeters -> first for horinzontal and second for vertical"
#    -> bold the whole paragraph (pPr/rPr + every run's rPr)
# ---------------------------------------------------------------------------
$idx = Get-ParaIndex "paremeters -> first for"
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">2 </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>paremeters</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve"> -&gt; first for </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>horinzontal</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve"> and second for vertical</w:t>
  </w:r>
</w:p>
'@
$null = $d.Paragraphs($idx).Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 4) "4 parameters -> start, top, end, bottom"
#    -> bold the whole paragraph (pPr/rPr + run rPr)
# ---------------------------------------------------------------------------
$idx = Get-ParaIndex "4 parameters -> start, top, end, bottom"
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>4 parameters -&gt; start, top, end, bottom</w:t>
  </w:r>
</w:p>
'@
$null = $d.Paragraphs($idx).Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 5) Insert a new bold "Note:" paragraph right before
#    "In jetpack, we use paddings for margin"
# ---------------------------------------------------------------------------
$idx = Get-ParaIndex "In jetpack, we use paddings for margin"
$null = $d.Paragraphs($idx).Range.InsertParagraphBefore()
$noteXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Note:</w:t>
  </w:r>
</w:p>
'@
$idx = Get-ParaIndex "In jetpack, we use paddings for margin"
$null = $d.Paragraphs($idx - 1).Range.InsertXML($noteXml)

# ---------------------------------------------------------------------------
# 6) Insert a new ".border" paragraph right after
#    "In jetpack, we use paddings for margin" (unchanged paragraph itself)
# ---------------------------------------------------------------------------
$idx = Get-ParaIndex "In jetpack, we use paddings for margin"
$null = $d.Paragraphs($idx).Range.InsertParagraphAfter()
$borderXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>.border</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
</w:p>
'@
$null = $d.Paragraphs($idx + 1).Range.InsertXML($borderXml)

Write-Output "edit complete"
